$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 467: corrected daily new-case count (2 -> 13); cumulative formula in B467 recalculates
$ws.Range("C467").Value = 13

# Row 468: corrected daily new-case count (1 -> 10); cumulative formula in B468 recalculates
$ws.Range("C468").Value = 10

# Row 469: fill in the day's data (previously blank placeholder row)
$ws.Range("C469").Value = 1

# L469/M469 are formatted as Text ("@") but need a genuine numeric 0 stored
# (matching how the rest of the column is populated) rather than the literal
# string "0" that a plain .Value assignment would produce on a Text cell.
# Temporarily switch to General, write the number, then restore the Text format.
$ws.Range("L469").NumberFormat = "general"
$ws.Range("L469").Value = 0
$ws.Range("L469").NumberFormat = "@"

$ws.Range("M469").NumberFormat = "general"
$ws.Range("M469").Value = 0
$ws.Range("M469").NumberFormat = "@"

# Update the saved selection on the frozen bottom-right pane
$ws.Range("A2").Select()
